$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01513913556067221
$ws.Range("C2").Value = 0.1935147000787938

$ws.Range("B3").Value = 0.07780334607687329
$ws.Range("C3").Value = 0.1377434218558929

$ws.Range("B4").Value = 0.5903019223163911
$ws.Range("C4").Value = 0.1704609998053457

$ws.Range("B5").Value = 0.9642806534098475
$ws.Range("C5").Value = 0.4099100892513007

$ws.Range("B6").Value = 0.8728958730672948
$ws.Range("C6").Value = 0.6585849000275269

$ws.Range("B7").Value = 0.6477492993357922
$ws.Range("C7").Value = 0.1121778859341113

$ws.Range("B8").Value = 0.004802521467208862
$ws.Range("C8").Value = 0.1970614433288574
